# Lesson 06 EJB deck — add the new "Transaction Handling" slide right
# before the closing "Git Repository Modules" slide (commit: "slide on
# transaction handling").

$p = $ppt.ActivePresentation

# The new slide re-uses the same "Title and Content" layout as the other
# content slides (e.g. the "Transactions" slide, slide 6).
$layout = $p.Slides.Item(6).CustomLayout

# Insert it at position 7 -> pushes the old slide 7 ("Git Repository
# Modules") down to position 8.
$newSlide = $p.Slides.AddSlide(7, $layout)

# --- Title -----------------------------------------------------------
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Transaction Handling"
$titleShape.TextFrame.TextRange.LanguageID = "en-US"

# --- Body content ------------------------------------------------------
$bodyShape = $newSlide.Shapes.Item(2)

# Custom placeholder position/size (matches the authored slide).
$bodyShape.Left = 16.578976377952756
$bodyShape.Top = 143.74992125984252
$bodyShape.Width = 928.8947244094488
$bodyShape.Height = 386.30259842519683

$tr = $bodyShape.TextFrame.TextRange

$lines = @(
  "EJB public methods will handle transactions by default",
  "Can use @annotations to fine tune them",
  "REQUIRED: default setting, start new transaction if none is active, or join current active one",
  "SUPPORTS: if there is an ongoing transaction, join it",
  "REQUIRES_NEW: always start a new transaction. If any ongoing, suspend them first",
  "MANDATORY: must be run in an ongoing transaction, otherwise fail",
  "NOT_SUPPORTED: put any ongoing transaction on hold",
  "NEVER: throw exception if in a transaction"
)
$tr.Text = [string]::Join("`r", $lines)
$tr.LanguageID = "en-US"

# Sub-bullets (REQUIRED / SUPPORTS / REQUIRES_NEW / MANDATORY /
# NOT_SUPPORTED / NEVER) sit one indent level in, with the annotation
# keyword itself italicised.
$keywords = @{
  3 = "REQUIRED"
  4 = "SUPPORTS"
  5 = "REQUIRES_NEW"
  6 = "MANDATORY"
  7 = "NOT_SUPPORTED"
  8 = "NEVER"
}

foreach ($idx in 3..8) {
  $para = $tr.Paragraphs($idx)
  $para.IndentLevel = 2

  $kw = $keywords[$idx]
  $kwRange = $para.Characters(1, $kw.Length)
  $kwRange.Font.Italic = $true
}
